# Refresh the cryptos price/volume snapshot (GitHub Actions cron update),
# including the Celestia/WEMIXToken row swap at rows 34-35.
# Price-column values that look numeric ("1.00", "306.48", ...) are forced
# to stay text (matching the source inlineStr cells) by toggling the
# NumberFormat to "@" before the write and clearing it right after, so no
# stray cell style lingers afterwards.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.800.01"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.333.80"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.48"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.509"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -5.39%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.511"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.96"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.17"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0799"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("E14").Value = "  -3.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.61"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.69%  "
$ws.Range("D16").Value = "2.308.18"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.796"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "42.734.40"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.61"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -7.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.69"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.24"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.89"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("E28").Value = "  +6.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.35"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.56"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.11"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.43%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.39"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.45"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0728"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.59"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +5.83%  "
$ws.Range("E38").Value = "  -4.77%  "
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("E40").Value = "  -4.32%  "
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.34"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("D43").Value = "2.014.15"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.90"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.36"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.93"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.94"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.91"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "2.560.69"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("E51").Value = "  +2.68%  "
